# Update the "Cong viec da dat duoc" (work completed) cell for the
# 02/11/2025 row (row 4) with the newly written progress note, wrap the
# text so the multi-line content displays properly, grow the row to fit
# it, and leave the active selection on the edited cell - matching the
# author's edit described in the commit message:
# "Update Frontend , Refactor code, Add Sequence Diagram and Activity Diagram"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newNote = "1. Chỉnh sửa lại mô hình dữ liệu mức ý niệm, thể hiện rõ các thực thể và mối quan hệ trong hệ thống bán vé sự kiện, đảm bảo tính logic và phù hợp với quy trình nghiệp vụ.`n2.  Hoàn thiện các sơ đồ tuần tự cho các chức năng.`n3. Hoàn thiện các sơ đồ hoạt động cho các chức năng.`n4.  Tiếp tục phát triền giao diện FrontEnd và refactor lại code dự án."

$cell = $ws.Range("D4")
$cell.Value = $newNote
$cell.WrapText = $true

$ws.Rows.Item(4).RowHeight = 152

[void]$cell.Select()

Write-Output "Updated D4 with the new progress note."
